# Update BOC USD rates (auto)
# Append the newly-published rate row to "All Published Values" and roll
# the day's publish count forward on "Daily Summary".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("All Published Values")

$newRow = 23

# Force the new row to text ("@") before writing so values like "697.85"
# or "2026-01-03" land as literal text, matching every other data row on
# this sheet instead of being auto-converted to numbers/dates.
$ws.Range("A" + $newRow + ":J" + $newRow).NumberFormat = "@"

$ws.Range("A" + $newRow).Value = "2026-01-03"
$ws.Range("B" + $newRow).Value = "2026-01-03 10:30:00"
$ws.Range("C" + $newRow).Value = "697.85"
$ws.Range("D" + $newRow).Value = "697.85"
$ws.Range("E" + $newRow).Value = "700.79"
$ws.Range("F" + $newRow).Value = "700.79"
$ws.Range("G" + $newRow).Value = "702.88"
$ws.Range("H" + $newRow).Value = "2026/01/03 10:30:00"
$ws.Range("I" + $newRow).Value = "2026-01-03 02:57:39"
$ws.Range("J" + $newRow).Value = "https://www.bankofchina.com/sourcedb/whpj/enindex_1619.html"

# Drop the temporary text format again so the new cells stay unstyled,
# same as the existing data rows.
$ws.Range("A" + $newRow + ":J" + $newRow).ClearFormats()

# Re-apply AutoFilter over the now-larger range (A1:J23).
$ws.AutoFilterMode = $false
$ws.Range("A1:J" + $newRow).AutoFilter() | Out-Null

# AutoFilter() alone doesn't rewrite the workbook-level hidden
# _FilterDatabase defined name, so update it explicitly to track the
# new range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "All Published Values!_FilterDatabase") {
        $n.RefersTo = "='All Published Values'!`$A`$1:`$J`$" + $newRow
    }
}

# Daily Summary: bump 2026-01-03's publish count (3 -> 4) now that a
# fourth rate was captured for that day.
$ws2 = $wb.Worksheets.Item("Daily Summary")
$ws2.Range("B5").Value = 4
